$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/29/2025  Through  10/5/2025"

# --- Weekly crime table updates (rows 14-21, 24-30) ---
# Row 14
$ws.Range("G14").Value = 2
$ws.Range("J14").Value = 4

# Row 15
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 25
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = 150

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("D16").Value = 2
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 54
$ws.Range("J16").Value = 65
$ws.Range("K16").Value = -16.923076923076
$ws.Range("L16").Value = -18.181818181818

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -41.666666666666
$ws.Range("I17").Value = 199
$ws.Range("J17").Value = 181
$ws.Range("K17").Value = 9.944751381215
$ws.Range("L17").Value = 7.567567567567

# Row 18
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("C33").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 1
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 44
$ws.Range("J18").Value = 56
$ws.Range("K18").Value = -21.428571428571
$ws.Range("L18").Value = -25.423728813559

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 166.666666666667
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = 18.75
$ws.Range("I19").Value = 172
$ws.Range("J19").Value = 203
$ws.Range("K19").Value = -15.270935960591
$ws.Range("L19").Value = -28.033472803347

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -31.25
$ws.Range("I20").Value = 128
$ws.Range("J20").Value = 136
$ws.Range("K20").Value = -5.882352941176
$ws.Range("L20").Value = -20.496894409937

# Row 21
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 25
$ws.Range("F21").Value = 61
$ws.Range("G21").Value = 67
$ws.Range("H21").Value = -8.955223880597
$ws.Range("I21").Value = 622
$ws.Range("J21").Value = 655
$ws.Range("K21").Value = -5.038167938931
$ws.Range("L21").Value = -13.850415512465

# Row 24
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 39
$ws.Range("G24").Value = 39
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 446
$ws.Range("J24").Value = 426
$ws.Range("K24").Value = 4.694835680751
$ws.Range("L24").Value = -4.496788008565

# Row 25
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -75
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = -27.272727272727
$ws.Range("I25").Value = 80
$ws.Range("J25").Value = 121
$ws.Range("K25").Value = -33.884297520661
$ws.Range("L25").Value = -37.984496124031

# Row 26
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 150
$ws.Range("F26").Value = 27
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 353
$ws.Range("J26").Value = 339
$ws.Range("K26").Value = 4.129793510324
$ws.Range("L26").Value = 19.256756756756

# Row 27
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 31
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = 63.157894736842
$ws.Range("L27").Value = 82.35294117647

# Row 28
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C33").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -66.666666666666
$ws.Range("I28").Value = 14
$ws.Range("J28").Value = 19
$ws.Range("K28").Value = -26.315789473684
$ws.Range("L28").Value = -33.333333333333

# Row 29
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$ws.Range("C33").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = -100
$ws.Range("I29").Value = 2
$ws.Range("J29").Value = 14
$ws.Range("K29").Value = -85.714285714285
$ws.Range("L29").Value = -60

# Row 30
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "0"
$ws.Range("C33").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -100
$ws.Range("I30").Value = 2
$ws.Range("J30").Value = 9
$ws.Range("K30").Value = -77.777777777777
$ws.Range("L30").Value = -50

$excel.CutCopyMode = $false